$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value2 = 59        # Total Trades: 58 -> 59
$wsSummary.Range("B9").Value2 = 49.15     # Win Rate %:   50 -> 49.15

# ---------------------------------------------------------------------
# Sheet: Strategy Status
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D5").Value2 = 26         # MarketMaking Trades:   25 -> 26
$wsStatus.Range("G5").Value2 = 57.69      # MarketMaking Win Rate: 60 -> 57.69

# ---------------------------------------------------------------------
# Sheet: All Trades
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Update existing row 60 (trade closed)
$wsAll.Range("G60").Value2 = 0.3
$wsAll.Range("H60").Value2 = "CLOSED"
$wsAll.Range("K60").Value2 = 100.59
$wsAll.Range("L60").Value2 = "early_exit"
$wsAll.Range("M60").Value2 = 0.14

# Append new row 93 for the newly opened trade
$wsAll.Range("A93").Value2 = 92
$wsAll.Range("B93").Value2 = "'2026-02-17"
$wsAll.Range("B93").Style = "Normal"
$wsAll.Range("C93").Value2 = "'20:57:13"
$wsAll.Range("C93").Style = "Normal"
$wsAll.Range("D93").Value2 = "MarketMaking"
$wsAll.Range("E93").Value2 = "UP"
$wsAll.Range("F93").Value2 = 0.3
$wsAll.Range("H93").Value2 = "OPEN"
$wsAll.Range("I93").Value2 = 0
$wsAll.Range("J93").Value2 = 0
$wsAll.Range("K93").Value2 = 100.5855022889912
$wsAll.Range("M93").Value2 = 0
$wsAll.Range("N93").Value2 = 0
$wsAll.Range("O93").Value2 = 0
$wsAll.Range("P93").Value2 = 0.6
$wsAll.Range("Q93").Value2 = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# Sheet: MarketMaking
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Update existing row 27 (trade closed)
$wsMM.Range("G27").Value2 = 0.3
$wsMM.Range("H27").Value2 = "CLOSED"
$wsMM.Range("K27").Value2 = 100.59
$wsMM.Range("P27").Value2 = "early_exit"
$wsMM.Range("Q27").Value2 = 0.14

# Append new row 60 for the newly opened trade
$wsMM.Range("A60").Value2 = 92
$wsMM.Range("B60").Value2 = "'2026-02-17"
$wsMM.Range("B60").Style = "Normal"
$wsMM.Range("C60").Value2 = "'20:57:13"
$wsMM.Range("C60").Style = "Normal"
$wsMM.Range("D60").Value2 = "MarketMaking"
$wsMM.Range("E60").Value2 = "UP"
$wsMM.Range("F60").Value2 = 0.3
$wsMM.Range("H60").Value2 = "OPEN"
$wsMM.Range("I60").Value2 = 0
$wsMM.Range("J60").Value2 = 0
$wsMM.Range("K60").Value2 = 100.5855022889912
$wsMM.Range("L60").Value2 = 0
$wsMM.Range("M60").Value2 = 0
$wsMM.Range("N60").Value2 = 0.6
$wsMM.Range("O60").Value2 = "Normal spread capture: 19600 bps"
$wsMM.Range("Q60").Value2 = 0
